$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 77388.66797673712
$ws.Range("D2").Value = 9992.97670278544
$ws.Range("E2").Value = 1770
$ws.Range("F2").Value = 25654.70836183914

$ws = $wb.Worksheets.Item("Capacities")
$ws.Range("C3").Value = 68

$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G2").Value = 13.6
$ws.Range("H2").Value = 27.2
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = 40.8
$ws.Range("K2").Value = 47.6
$ws.Range("L2").Value = 54.4
$ws.Range("M2").Value = 61.2
$ws.Range("N2").Value = 68
$ws.Range("O2").Value = 61.2
$ws.Range("P2").Value = 54.4
$ws.Range("Q2").Value = 47.6
$ws.Range("R2").Value = 34
$ws.Range("S2").Value = 20.4
$ws.Range("T2").Value = 13.6
$ws.Range("I3").Value = 27.2
$ws.Range("J3").Value = 40.8
$ws.Range("K3").Value = 54.4
$ws.Range("L3").Value = 61.2
$ws.Range("M3").Value = 68
$ws.Range("N3").Value = 54.4
$ws.Range("O3").Value = 47.6
$ws.Range("P3").Value = 34
$ws.Range("Q3").Value = 34
$ws.Range("R3").Value = 20.4
$ws.Range("S3").Value = 13.6
$ws.Range("J4").Value = 6.8
$ws.Range("K4").Value = 27.2
$ws.Range("L4").Value = 47.6
$ws.Range("M4").Value = 54.4
$ws.Range("N4").Value = 37.5831241710018
$ws.Range("O4").Value = 47.6
$ws.Range("P4").Value = 27.2
$ws.Range("Q4").Value = 13.6
$ws.Range("R4").Value = 6.8

$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("G2").Value = 64.3
$ws.Range("H2").Value = 14.2
$ws.Range("I2").Value = 2.8
$ws.Range("J2").Value = 1.8
$ws.Range("K2").Value = 21.6
$ws.Range("L2").Value = 33.6
$ws.Range("M2").Value = 37.8
$ws.Range("N2").Value = 42
$ws.Range("O2").Value = 30
$ws.Range("P2").Value = 136.2342720130611
$ws.Range("Q2").Value = 21.6
$ws.Range("R2").Value = 0.2
$ws.Range("S2").Value = 32.4
$ws.Range("T2").Value = 45.6
$ws.Range("I3").Value = 27.2
$ws.Range("J3").Value = 40.8
$ws.Range("K3").Value = 54.4
$ws.Range("L3").Value = 61.2
$ws.Range("M3").Value = 44.6
$ws.Range("N3").Value = 28.4
$ws.Range("O3").Value = 47.6
$ws.Range("P3").Value = 5.630792776247632
$ws.Range("Q3").Value = 8
$ws.Range("R3").Value = 20.4
$ws.Range("J4").Value = 6.8
$ws.Range("K4").Value = 27.2
$ws.Range("L4").Value = 47.6
$ws.Range("M4").Value = 31
$ws.Range("N4").Value = 37.5831241710018
$ws.Range("O4").Value = 47.6
$ws.Range("P4").Value = 27.2
$ws.Range("Q4").Value = 13.6
$ws.Range("R4").Value = 6.8

$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("G2").Value = 183.657
$ws.Range("H2").Value = 197.715
$ws.Range("I2").Value = 200.487
$ws.Range("J2").Value = 202.269
$ws.Range("K2").Value = 223.653
$ws.Range("L2").Value = 256.917
$ws.Range("M2").Value = 294.3389999999999
$ws.Range("N2").Value = 335.9189999999999
$ws.Range("O2").Value = 365.6189999999999
$ws.Range("P2").Value = 500.4909292929304
$ws.Range("Q2").Value = 521.8749292929305
$ws.Range("R2").Value = 522.0729292929304
$ws.Range("S2").Value = 554.1489292929305
$ws.Range("B3").Value = 172.5252525252518
$ws.Range("C3").Value = 152.8282828282828
$ws.Range("D3").Value = 139.6969696969689
$ws.Range("E3").Value = 139.6969696969689
$ws.Range("F3").Value = 139.6969696969689
$ws.Range("G3").Value = 120
$ws.Range("H3").Value = 120
$ws.Range("I3").Value = 146.928
$ws.Range("J3").Value = 187.32
$ws.Range("K3").Value = 241.176
$ws.Range("L3").Value = 301.764
$ws.Range("M3").Value = 345.918
$ws.Range("N3").Value = 374.034
$ws.Range("O3").Value = 421.158
$ws.Range("P3").Value = 426.7324848484852
$ws.Range("Q3").Value = 434.6524848484852
$ws.Range("R3").Value = 454.8484848484852
$ws.Range("S3").Value = 454.8484848484852
$ws.Range("T3").Value = 323.5353535353539
$ws.Range("U3").Value = 323.5353535353539
$ws.Range("V3").Value = 323.5353535353539
$ws.Range("W3").Value = 244.7474747474744
$ws.Range("X3").Value = 244.7474747474744
$ws.Range("Y3").Value = 205.3535353535346
$ws.Range("J4").Value = 126.732
$ws.Range("K4").Value = 153.66
$ws.Range("L4").Value = 200.784
$ws.Range("M4").Value = 231.474
$ws.Range("N4").Value = 268.6812929292918
$ws.Range("O4").Value = 315.8052929292918
$ws.Range("P4").Value = 342.7332929292918
$ws.Range("Q4").Value = 356.1972929292918

$ws = $wb.Worksheets.Item("Feed in from Type 2")
$ws.Range("J2").Value = 0
$ws.Range("P2").Value = 110.4342720130611
$ws.Range("N3").Value = 0
$ws.Range("P3").Value = 0.2307927762476334
$ws.Range("P4").Value = 0

$ws = $wb.Worksheets.Item("Feed in from Type 3")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0

$ws = $wb.Worksheets.Item("Feed in from Type 4")
$ws.Range("J2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("S2").Value = 53.6
$ws.Range("J3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("S3").Value = 9.6
$ws.Range("J4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
